# Scheduled runner refresh: update market-price-derived Leve profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) with latest data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 97
$ws.Range("H97").Value = 1984.3636
$ws.Range("J97").Value = 2082.8
$ws.Range("L97").Value = 6248.400000000001
$ws.Range("N97").Value = -7240.400000000001

# Row 99
$ws.Range("H99").Value = 1137
$ws.Range("I99").Value = 520
$ws.Range("J99").Value = 1959.6666
$ws.Range("K99").Value = 1560
$ws.Range("L99").Value = 5878.9998
$ws.Range("M99").Value = -62
$ws.Range("N99").Value = -8874.9998

# Row 138
$ws.Range("H138").Value = 5001610.5
$ws.Range("I138").Value = 1658.1578
$ws.Range("J138").Value = 9525377
$ws.Range("K138").Value = 4974.4734
$ws.Range("L138").Value = 28576131
$ws.Range("M138").Value = 165.5266000000001
$ws.Range("N138").Value = -28586411

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5316.96
$ws.Range("I32").Value = 2766.4768
$ws.Range("J32").Value = 20984.215
$ws.Range("K32").Value = 2766.4768
$ws.Range("L32").Value = 20984.215
$ws.Range("M32").Value = -2479.4768
$ws.Range("N32").Value = -21558.215

# Row 45
$ws.Range("H45").Value = 1049.9
$ws.Range("I45").Value = 971.4286
$ws.Range("J45").Value = 1233
$ws.Range("K45").Value = 971.4286
$ws.Range("L45").Value = 1233
$ws.Range("M45").Value = -594.4286
$ws.Range("N45").Value = -1987

# Row 61
$ws.Range("H61").Value = 30365458
$ws.Range("I61").Value = 38501320
$ws.Range("J61").Value = 146544.86
$ws.Range("K61").Value = 38501320
$ws.Range("L61").Value = 146544.86
$ws.Range("M61").Value = -38501108
$ws.Range("N61").Value = -146968.86

# Row 110
$ws.Range("H110").Value = 1739.2
$ws.Range("I110").Value = 1300
$ws.Range("J110").Value = 1849
$ws.Range("K110").Value = 1300
$ws.Range("L110").Value = 1849
$ws.Range("M110").Value = 745
$ws.Range("N110").Value = -5939

# Row 122
$ws.Range("H122").Value = 6537664
$ws.Range("I122").Value = 1605.1111
$ws.Range("J122").Value = 13890731
$ws.Range("K122").Value = 4815.3333
$ws.Range("L122").Value = 41672193
$ws.Range("M122").Value = -2365.3333
$ws.Range("N122").Value = -41677093

# Row 132
$ws.Range("H132").Value = 92909.22
$ws.Range("I132").Value = 60465.53
$ws.Range("J132").Value = 184833
$ws.Range("K132").Value = 181396.59
$ws.Range("L132").Value = 554499
$ws.Range("M132").Value = -178866.59
$ws.Range("N132").Value = -559559

# Row 136
$ws.Range("H136").Value = 30365458
$ws.Range("I136").Value = 38501320
$ws.Range("J136").Value = 146544.86
$ws.Range("K136").Value = 115503960
$ws.Range("L136").Value = 439634.58
$ws.Range("M136").Value = -115501410
$ws.Range("N136").Value = -444734.58

$ws = $wb.Worksheets.Item("BSM")
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("N132").Value = 0

# Row 133
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120

# Row 134
$ws.Range("H134").Value = 2223.4375
$ws.Range("I134").Value = 2223.4375
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6670.3125
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -4135.3125

# Row 137
$ws.Range("H137").Value = 54225
$ws.Range("J137").Value = 54225
$ws.Range("L137").Value = 54225
$ws.Range("N137").Value = -64425

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3089.9688
$ws.Range("I31").Value = 2612.7222
$ws.Range("J31").Value = 3703.5715
$ws.Range("K31").Value = 2612.7222
$ws.Range("L31").Value = 3703.5715
$ws.Range("M31").Value = -2317.7222
$ws.Range("N31").Value = -4293.5715

# Row 34
$ws.Range("H34").Value = 3089.9688
$ws.Range("I34").Value = 2612.7222
$ws.Range("J34").Value = 3703.5715
$ws.Range("K34").Value = 2612.7222
$ws.Range("L34").Value = 3703.5715
$ws.Range("M34").Value = -2410.7222
$ws.Range("N34").Value = -4107.5715

# Row 58
$ws.Range("H58").Value = 45456316
$ws.Range("I58").Value = 45456316
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 45456316
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -45456113

# Row 107
$ws.Range("H107").Value = 453.81818
$ws.Range("I107").Value = 428.42105
$ws.Range("J107").Value = 614.6667
$ws.Range("K107").Value = 428.42105
$ws.Range("L107").Value = 614.6667
$ws.Range("M107").Value = 1491.57895
$ws.Range("N107").Value = -4454.6667

# Row 136
$ws.Range("H136").Value = 45456316
$ws.Range("I136").Value = 45456316
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 136368948
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -136366398
$ws.Range("N136").Value = -136366398

# Row 140
$ws.Range("H140").Value = 55700
$ws.Range("J140").Value = 55700
$ws.Range("L140").Value = 55700
$ws.Range("N140").Value = -66060

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1647.875
$ws.Range("I102").Value = 1599.8572
$ws.Range("J102").Value = 1984
$ws.Range("K102").Value = 1599.8572
$ws.Range("L102").Value = 1984
$ws.Range("M102").Value = 22.14280000000008
$ws.Range("N102").Value = -5228

# Row 122
$ws.Range("H122").Value = 2064.3333
$ws.Range("I122").Value = 1915.7
$ws.Range("K122").Value = 5747.1
$ws.Range("M122").Value = -3297.1

# Row 132
$ws.Range("H132").Value = 59823.83
$ws.Range("I132").Value = 37097.31
$ws.Range("J132").Value = 169668.67
$ws.Range("K132").Value = 111291.93
$ws.Range("L132").Value = 509006.01
$ws.Range("M132").Value = -108761.93
$ws.Range("N132").Value = -514066.01

# Row 136
$ws.Range("H136").Value = 39054.332
$ws.Range("J136").Value = 39054.332
$ws.Range("L136").Value = 117162.996
$ws.Range("N136").Value = -122262.996

$ws = $wb.Worksheets.Item("LTW")
# Row 53
$ws.Range("H53").Value = 9500
$ws.Range("J53").Value = 9500
$ws.Range("L53").Value = 9500
$ws.Range("N53").Value = -10536

# Row 132
$ws.Range("H132").Value = 49890.715
$ws.Range("I132").Value = 1637.9375
$ws.Range("J132").Value = 204299.6
$ws.Range("K132").Value = 4913.8125
$ws.Range("L132").Value = 612898.8
$ws.Range("M132").Value = -2383.8125
$ws.Range("N132").Value = -617958.8

# Row 135
$ws.Range("H135").Value = 49200
$ws.Range("J135").Value = 49200
$ws.Range("L135").Value = 49200
$ws.Range("N135").Value = -59340

# Row 136
$ws.Range("H136").Value = 144930.36
$ws.Range("I136").Value = 101301.5
$ws.Range("J136").Value = 254002.5
$ws.Range("K136").Value = 303904.5
$ws.Range("L136").Value = 762007.5
$ws.Range("M136").Value = -301354.5
$ws.Range("N136").Value = -767107.5

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1884
$ws.Range("I126").Value = 1845.7778
$ws.Range("K126").Value = 5537.3334
$ws.Range("M126").Value = -3067.3334

# Row 136
$ws.Range("H136").Value = 64662.125
$ws.Range("I136").Value = 48641.24
$ws.Range("J136").Value = 95247.45
$ws.Range("K136").Value = 145923.72
$ws.Range("L136").Value = 285742.35
$ws.Range("M136").Value = -143373.72
$ws.Range("N136").Value = -290842.35
